$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Update the conversion factor and its accompanying label.
$ws.Range("A26").Value = 0.75350342301658668
$ws.Range("B26").Value = "2023 dollars per 2012 dollar"

# Move the active selection from B30 to A26, matching the saved sheet view.
$ws.Range("A26").Select()
